$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Median Salary (City)" header in F1, matching the style of the
# other header cells (bold font, border, centered alignment) by copying the
# formatting from the adjacent E1 header cell.
$ws.Range("F1").Value = "Median Salary (City)"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column E ("Median Salary (National)") switches from numeric values to text
# values holding the same digits. A leading apostrophe forces Excel to store
# the entry as text instead of re-parsing it back into a number.
$ws.Range("E2:E11").Value = "'103500"
$ws.Range("E2:E11").Style = "Normal"

# New column F data cells (F2:F11) are empty text cells (present in the
# sheet, but holding no content). A lone leading apostrophe creates an empty
# text cell without leaving stray content behind.
$ws.Range("F2:F11").Value = "'"
$ws.Range("F2:F11").Style = "Normal"
